# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows that changed after repulling data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 10
$ws.Range("F7").Value = -5
$ws.Range("F9").Value = -11
$ws.Range("F12").Value = -3
$ws.Range("F14").Value = 10
$ws.Range("F15").Value = 9
$ws.Range("F16").Value = 10
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = -3
$ws.Range("F21").Value = -1
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = 2
